# "Generate Report for Archive"
#
# The localization status for the two tracked files moved on from
# "Ready for handoff" to "In Translation". Update the Status text on the
# per-language sheets (which also drives the Overview roll-up, since both
# point at the same shared string), then let the report's column-fit logic
# re-run so the now-shorter status text no longer needs such wide columns.

$wb = $excel.ActiveWorkbook

$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Status column ("C") on each per-language sheet.
$zhSheet.Range("C2:C3").Value = "In Translation"
$deSheet.Range("C2:C3").Value = "In Translation"

# The Overview sheet mirrors the same status text in its zh-cn/de-de
# columns ("E"/"F").
$overview.Range("E2:F3").Value = "In Translation"

# Re-fit the Status column(s) now that the text is shorter.
$zhSheet.Columns.Item(3).ColumnWidth = 12.5
$deSheet.Columns.Item(3).ColumnWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
